$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.730.93"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "2.263.85"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.71%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0797"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("D15").Value = "2.614.21"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").Value = "2.239.04"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.765"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.68%  "

$ws.Range("D19").Value = "41.641.02"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.05%  "

$ws.Range("D21").Value = "0.0₃0904"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  +2.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.34%  "

$ws.Range("E30").Value = "  -4.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.03%  "

$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("E35").Value = "  +1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("E40").Value = "  -0.76%  "

$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("E42").Value = "  +0.50%  "

$ws.Range("D43").Value = "2.026.37"
$ws.Range("E43").Value = "  -3.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.93%  "

$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("E47").Value = "  +12.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.75%  "

$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
